# Added 4wk low sales check: update Seasonality Index (column L) values
# on the "Forecast Comparison" sheet to reflect recalculated figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

$ws.Range("L2").Value  = 1.15
$ws.Range("L3").Value  = 1.05
$ws.Range("L5").Value  = 1.01
$ws.Range("L6").Value  = 1.13
$ws.Range("L7").Value  = 1.19
$ws.Range("L8").Value  = 1.01
$ws.Range("L9").Value  = 0.88
$ws.Range("L10").Value = 0.86
$ws.Range("L11").Value = 1.19
$ws.Range("L12").Value = 0.88
$ws.Range("L13").Value = 0.94
$ws.Range("L14").Value = 1.01
$ws.Range("L15").Value = 0.89
$ws.Range("L16").Value = 1.08
$ws.Range("L17").Value = 1.07
